# Revert "Revert "Merge pull request #48 from LakeFishing/main""
#
# Changes applied to the active worksheet:
#   - F4, B9, B14 previously held the shared string "time"; they are
#     updated to hold the string "special".
#   - The saved cursor/selection moves from the old multi-area
#     selection (C7 / C12) to a single-cell selection on B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "special"
$ws.Range("B9").Value = "special"
$ws.Range("B14").Value = "special"

[void]$ws.Range("B14").Select()
